# Adding Master Data XLS
# Appends the French-localized rows (id 10013-10018) to the
# master-app_detail sheet, mirroring the existing English/Arabic rows,
# and restores the column widths / selection / page setup that Excel
# wrote when the data was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the six new French rows (rows 14-19) -----------------------
$newRows = @(
  @(10013, "Pré-inscription",      "Portail Web pour les pré-inscriptions",                                 "fra", $true, "superadmin", "now()"),
  @(10014, "Client dinscription",  "Application de bureau pour les inscriptions",                            "fra", $true, "superadmin", "now()"),
  @(10015, "Processeur dinscription", "Demande de post-inscription",                                         "fra", $true, "superadmin", "now()"),
  @(10016, "Authentification ID",  "Application pour lauthentification du fournisseur de services tiers",   "fra", $true, "superadmin", "now()"),
  @(10017, "Contrôle didentité",   "Portail Web pour la configuration dapplications",                        "fra", $true, "superadmin", "now()"),
  @(10018, "Portail Résident",     "Portail Web pour les services de génération de post-ID",                 "fra", $true, "superadmin", "now()")
)

$r = 14
foreach ($row in $newRows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $r++
}

# --- Column widths (id column best-fit, name column widened) -----------
# The engine snaps ColumnWidth to whole-pixel steps, same as Excel; these
# inputs land on the closest achievable width to the author's saved
# 5.81640625 / 19.26953125 character widths.
$ws.Columns.Item(1).ColumnWidth = 5.0
$ws.Columns.Item(2).ColumnWidth = 18.5

# --- Scroll / selection state as left by the author on save ------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A20:XFD1048576").Select()

# --- Page setup ----------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

$wb.Save()
